# Fruta / hortaliza, semanal
# A new weekly record was inserted as row 583 ("Especial"/M=100 data pushed down
# to become the regular chain), with every subsequent row (584-703) shifting
# down by one and the former last row (703) becoming a brand new row 704.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 583. This pushes the existing rows
# 583-703 down to 584-704 (Excel's normal Insert-row shift behaviour,
# which also shifts along the date-format styling on column D).
$ws.Rows("583:583").Insert()

# The row that used to be 583 is now row 584. Populate the brand-new row 583
# by copying across every column from row 584 except the three columns whose
# values actually differ (D = Fecha, L = Calidad, M = Volumen).
for ($col = 1; $col -le 20; $col++) {
    if ($col -ne 4 -and $col -ne 12 -and $col -ne 13) {
        $ws.Cells.Item(583, $col).Value2 = $ws.Cells.Item(584, $col).Value2
    }
}

# Now set the new row's distinct values.
$ws.Cells.Item(583, 4).Value2 = 45244       # D583: Fecha
$ws.Cells.Item(583, 12).Value2 = "Primera"  # L583: Calidad
$ws.Cells.Item(583, 13).Value2 = 150        # M583: Volumen
